# Fill out the "Method Inputs" and "Expected Result" columns (E:G) for the
# unit-test-plan rows (7-16) of the Pixell client test plan.
#
# NOTE: the order of these assignments is deliberately chosen to reproduce
# the original author's shared-string insertion order (new distinct text
# values are interned in the order they were first typed into cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 27: "None" - used for every Preconditions cell E7:E16
$ws.Range("E7").Value = "None"

# 28: Row 7 Method Inputs
$ws.Range("F7").Value = 'Client_number=12345 first_name="Wendy"        last_name="ways" email_address="WendyWays@pixell-river.com"'

# 29: Row 8 Method Inputs
$ws.Range("F8").Value = 'Client_number="hello world" first_name="Wendy"        last_name="ways" email_address="WendyWays@pixell-river.com"'

# 30: Row 7 Expected Result
$ws.Range("G7").Value = "Attributes assigned"

# 31: Row 8 Expected Result - also reused by rows 9 & 10
$ws.Range("G8").Value = "ValueError is raised, error message printed."

# 32: Row 9 Method Inputs
$ws.Range("F9").Value = 'Client_number=12345 first_name=" "        last_name="ways" email_address="WendyWays@pixell-river.com"'

# 33: Row 10 Method Inputs
$ws.Range("F10").Value = 'Client_number=12345 first_name="Wendy"        last_name=" " email_address="WendyWays@pixell-river.com"'

# 34: Row 11 Method Inputs
$ws.Range("F11").Value = 'Client_number=12345 first_name="Wendy"        last_name=" " email_address="WendyWayspixell-river.com"'

# 35: Row 11 Expected Result
$ws.Range("G11").Value = 'Error is raised, "email@pixell-river.com" printed as default attribute'

# 36: Row 12 Method Inputs
$ws.Range("F12").Value = "Client_number=12345"

# 37: Row 12 Expected Result - also reused by rows 13, 14, 15 & 16
$ws.Range("G12").Value = "no error"

# 38: Row 13 Method Inputs
$ws.Range("F13").Value = 'first_name="Wendy"'

# 39: Row 14 Method Inputs
$ws.Range("F14").Value = 'last_name="Ways"'

# 40: Row 15 Method Inputs
$ws.Range("F15").Value = 'email_address="WendyWays@pixell-river.com"'

# 41: Row 16 Method Inputs
$ws.Range("F16").Value = "Ways, Wendy [12345] - WendyWays@pixell-river.com"

# Fill in the remaining "None" (Preconditions) and "no error"/"ValueError..."
# (Expected Result) cells that reuse already-interned shared strings.
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"
$ws.Range("E12").Value = "None"
$ws.Range("E13").Value = "None"
$ws.Range("E14").Value = "None"
$ws.Range("E15").Value = "None"
$ws.Range("E16").Value = "None"

$ws.Range("G9").Value = "ValueError is raised, error message printed."
$ws.Range("G10").Value = "ValueError is raised, error message printed."

$ws.Range("G13").Value = "no error"
$ws.Range("G14").Value = "no error"
$ws.Range("G15").Value = "no error"
$ws.Range("G16").Value = "no error"

# Leave the view focused where the author last left it (matches the saved
# sheetView selection in the workbook).
$ws.Activate()
$ws.Range("G16").Select()
